$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "About" sheet: refresh sourcing/notes block
# ---------------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

# Drop the old "last edited" date stamp that lived in column C and remove the
# now-unused column entirely (the refreshed layout only uses columns A:B).
$about.Columns("C").Delete()

# Make room for four new source-citation rows (rows 4-7) and three new
# note rows (rows 16-18, after the shift caused by the first insert).
$about.Rows("4:7").Insert()
$about.Rows("16:18").Insert()

# Rows 3-7: source + citation detail lines (written in this order so the
# shared-string table comes out in the same sequence the authored file uses).
$about.Range("B7").Value = "Table 5 Generalized Cost Coefficient Calibration"
$about.Range("B3").Value = "United States EPA"
$about.Range("B5").Value = "Consumer Vehicle Choice Model Documentation"
$about.Range("B6").Value = "https://nepis.epa.gov/Exe/ZyPDF.cgi/P100EZ37.PDF?Dockey=P100EZ37.PDF"
$about.Range("B4").Value = 2012
# The publication year is a plain left-aligned number, not the old date stamp
# that used to occupy this style slot.
$about.Range("B4").HorizontalAlignment = -4131

# Rows 16-17: new explanatory notes about the -3 / -5 logit exponent choice.
$about.Range("A16").Value = "We choose a value of -3 for passenger vehicles and a value of -5 for other vehicle types, "
$about.Range("A17").Value = "based on the ranges in Table 5 of the cited EPA documentation."

# ---------------------------------------------------------------------------
# "TTLE" sheet: logit exponents move from -3 to -5 for every row
# ---------------------------------------------------------------------------
$ttle = $wb.Worksheets.Item("TTLE")
$ttle.Range("B2:C7").Value = -5
